# Apply data_dictionary.xlsx normalization edits
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("records_2022")

# E5 ("value" row) - Missing value meaning: N/A -> Invalid entry
$ws.Range("E5").Value = "Invalid entry"

# E8 ("status" row) - Missing value meaning: N/A -> Gets defaulted to review
$ws.Range("E8").Value = "Gets defaulted to review"

# D7 ("source" row) - Allowed Values: drop "source_system" from the list
$ws.Range("D7").Value = "system_a, system_b, system_c, manual_entry, import_batch"

# E7 ("source" row) - Missing value meaning: N/A -> Gets defaulted to import_batch
$ws.Range("E7").Value = "Gets defaulted to import_batch"

# Make records_2022 the active sheet and move the selection to D4, matching
# the final cursor position left behind after the edits.
$ws.Activate()
$ws.Range("D4").Select() | Out-Null

$wb.Save()
